$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 363; this shifts the existing rows 363-393 down to 364-394
# (and their D/J/K/L/M/P values move with them, matching the rest of the diff).
$ws.Rows(363).Insert()

# Populate the newly inserted row 363 with the new weekly record.
$ws.Cells.Item(363, 1).Value = 10
$ws.Cells.Item(363, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(363, 3).Value = "La Araucanía"
$ws.Cells.Item(363, 4).Value = 45106
$ws.Cells.Item(363, 5).Value = 9
$ws.Cells.Item(363, 6).Value = 100112039
$ws.Cells.Item(363, 7).Value = "Ciboulette"
$ws.Cells.Item(363, 8).Value = "Sin especificar"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 90
$ws.Cells.Item(363, 11).Value = 6000
$ws.Cells.Item(363, 12).Value = 6000
$ws.Cells.Item(363, 13).Value = 6000
$ws.Cells.Item(363, 14).Value = "$/docena de atados"
$ws.Cells.Item(363, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(363, 16).Value = 2000
$ws.Cells.Item(363, 17).Value = 3
$ws.Cells.Item(363, 18).Value = "Hortaliza"
